$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (G:K) for the "meta" group, right before the existing
# "arrecadado_sucesso" group (old column G), shifting everything after it
# to the right by 5 columns.
$ws.Range("G1:K1").EntireColumn.Insert()

# Headers for the newly inserted columns
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# Match the currency formatting used by the neighbouring "arrecadado_*" columns
$ws.Range("G2:K2").NumberFormat = $ws.Range("L2").NumberFormat

# Data values for the new "meta" columns
$ws.Range("G2").Value = 13973042.60019265
$ws.Range("H2").Value = 16834.99108456945
$ws.Range("I2").Value = 17015.69760983049
$ws.Range("J2").Value = 31.89582864100442
$ws.Range("K2").Value = 189313.7035611726

# Recomputed stats shifted: the contribuicoes_std value changed very slightly
$ws.Range("W2").Value = 423.0192251466749
